$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old movie list data (rows 1-18, columns A-C)
$ws.Range("A1:C18").Clear()

# New header row: Title, Stars, Director, Plot
$ws.Range("A1").Value = "Title"
$ws.Range("B1").Value = "Stars"
$ws.Range("C1").Value = "Director"
$ws.Range("D1").Value = "Plot"

# Row 2 - November
$ws.Range("A2").Value = "November"
$ws.Range("B2").Formula = "=TEXT(7.8,""0.0"")"
$ws.Range("C2").Value = "Achero Mañas"
$ws.Range("D2").Value = "Impelled by a spirit which still preserves a patina of idealism, Alfredo arrives in Madrid intent on creating ""a performance which is freer, straight from the heart, capable of making ..."

# Row 3 - Outlaw King (info not found via omdbAPI)
$ws.Range("A3").Value = "Outlaw King"
$ws.Range("B3").Value = "Not found"
$ws.Range("C3").Value = "Not found"
$ws.Range("D3").Value = "Not found"

# Row 4 - The Nun
$ws.Range("A4").Value = "The Nun"
$ws.Range("B4").Formula = "=TEXT(4,""0.0"")"
$ws.Range("C4").Value = "Luis de la Madrid"
$ws.Range("D4").Value = "A group of teenage girls are terrorized by Sister Ursula, a nun that believes she must rid the world of all sin. After Sister Ursula mysteriously disappears, the Catholic school is shut ..."

# Convert the TEXT() formula results to plain text values (no formula, no number
# reinterpretation) so the Stars column stores "7.8"/"4.0" as text, matching
# how the ratings are stored elsewhere in the sheet.
$ws.Range("B2").Copy()
$ws.Range("B2").PasteSpecial(-4163)
$ws.Range("B4").Copy()
$ws.Range("B4").PasteSpecial(-4163)

$excel.CutCopyMode = 0
